$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set G29 to 1 (the "Particular effort" criteria gained 1 point)
$ws.Range("G29").Value = 1

# Add a comment / note in column J explaining the extra point
$ws.Range("J29").Value = "Use of semantic elements in index.html"

# Update the selected cell to reflect where the edit was made
$ws.Range("J26").Select()
